$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.782.69"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").Value = "'3.049.78"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.20%  "

$ws.Range("D5").Value = "'558.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'142.25"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.16%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "'3.046.81"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.17%  "

$ws.Range("D9").Value = "'0.513"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.55%  "

$ws.Range("D10").Value = "'0.153"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("D11").Value = "'6.18"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -12.19%  "

$ws.Range("D12").Value = "'0.493"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +6.76%  "

$ws.Range("D13").Value = "'0.0000229"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").Value = "'35.60"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("D15").Value = "'3.549.94"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("D16").Value = "'63.821.04"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.51%  "

$ws.Range("D17").Value = "'3.049.10"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.56%  "

$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").Value = "'6.77"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").Value = "'474.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("D21").Value = "'14.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.45%  "

$ws.Range("D22").Value = "'0.682"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").Value = "'14.52"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +9.54%  "

$ws.Range("D24").Value = "'7.54"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").Value = "'82.46"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.69%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.72%  "

$ws.Range("D27").Value = "'2.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("D28").Value = "'8.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.32%  "

$ws.Range("D29").Value = "'2.03"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("E30").Value = "  +0.41%  "

$ws.Range("D31").Value = "'26.22"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("E32").Value = "  -1.78%  "

$ws.Range("D33").Value = "'2.45"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").Value = "'5.74"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.84%  "

$ws.Range("D35").Value = "'6.21"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.92%  "

$ws.Range("D36").Value = "'54.48"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.08%  "

$ws.Range("D37").Value = "'0.0409"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").Value = "'446.87"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.78%  "

$ws.Range("D39").Value = "'0.0812"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("E40").Value = "  +3.67%  "

$ws.Range("D41").Value = "'3.009.17"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.117"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.26"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").Value = "'0.267"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.70%  "

$ws.Range("D45").Value = "'28.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.80%  "

$ws.Range("D46").Value = "'2.24"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +7.09%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").Value = "'0.113"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("D49").Value = "'117.76"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("D50").Value = "'0.0₃0513"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").Value = "'2.09"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.43%  "
